$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Battery_Data" ---
# Model now only reports a single "upgrade" (instead of up to 3), so the
# "... at upgrade 2" / "... at upgrade 3" rows for each of the three
# categories (Nominal Capacity, Investment, Yearly O&M Cost) are removed,
# keeping only the "... at upgrade 1" row for each. Delete bottom-up so
# the not-yet-deleted row numbers stay stable.
$ws1 = $wb.Worksheets.Item("Battery_Data")

$ws1.Rows.Item(13).EntireRow.Delete()   # Yearly O&M Cost at upgrade 3
$ws1.Rows.Item(12).EntireRow.Delete()   # Yearly O&M Cost at upgrade 2
$ws1.Rows.Item(10).EntireRow.Delete()   # Investment at upgrade 3
$ws1.Rows.Item(9).EntireRow.Delete()    # Investment at upgrade 2
$ws1.Rows.Item(7).EntireRow.Delete()    # Nominal Capacity at upgrade 3
$ws1.Rows.Item(6).EntireRow.Delete()    # Nominal Capacity at upgrade 2

# Fill in the real computed results (previously all placeholder zeros)
$ws1.Range("B5").Value = 6175.0338886374811
$ws1.Range("B6").Value = 2762.957163131955
$ws1.Range("B7").Value = 55.259143262639093
$ws1.Range("B8").Value = 144.6057673583835

# Re-fit the columns to the (now shorter) content
$ws1.Columns.Item(1).ColumnWidth = 36.333
$ws1.Columns.Item(2).ColumnWidth = 11.0

# --- Sheet 2: "Yearly BRC" ---
# Results now span 5 years instead of 3, with real computed values
# replacing the zero placeholders.
$ws2 = $wb.Worksheets.Item("Yearly BRC")

$ws2.Range("B2").Value = 72.376755136928352
$ws2.Range("B3").Value = 42.342319963364957
$ws2.Range("B4").Value = 21.39342939190227

# Add the two new year rows, copying the formatting of the last
# existing data row (row 4) so the new header cells match style "1"
$ws2.Range("A4").Copy() | Out-Null
$ws2.Range("A5:A6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws2.Range("A5").Value = "Battery Reposition Cost at y = 4"
$ws2.Range("B5").Value = 2.3731693227008011
$ws2.Range("A6").Value = "Battery Reposition Cost at y = 5"
$ws2.Range("B6").Value = 6.1200935434860089

# Re-fit the columns to the (now longer) content
$ws2.Columns.Item(1).ColumnWidth = 26.833
$ws2.Columns.Item(2).ColumnWidth = 11.0
